$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1277.909
$ws.Range("I40").Value = 776
$ws.Range("J40").Value = 1696.1666
$ws.Range("K40").Value = 776
$ws.Range("L40").Value = 1696.1666
$ws.Range("M40").Value = -601
$ws.Range("N40").Value = -2046.1666

$ws.Range("H64").Value = 4166.6665
$ws.Range("I64").Value = 3250
$ws.Range("J64").Value = 4625
$ws.Range("K64").Value = 3250
$ws.Range("L64").Value = 4625
$ws.Range("M64").Value = -3002
$ws.Range("N64").Value = -5121

$ws.Range("H67").Value = 4166.6665
$ws.Range("I67").Value = 3250
$ws.Range("J67").Value = 4625
$ws.Range("K67").Value = 3250
$ws.Range("L67").Value = 4625
$ws.Range("M67").Value = -2392
$ws.Range("N67").Value = -6341

$ws.Range("H74").Value = 5500
$ws.Range("J74").Value = 5500
$ws.Range("L74").Value = 5500
$ws.Range("N74").Value = -7372

$ws.Range("H77").Value = 5500
$ws.Range("J77").Value = 5500
$ws.Range("L77").Value = 27500
$ws.Range("N77").Value = -36860

$ws.Range("H80").Value = 12184832
$ws.Range("I80").Value = 2650.6667
$ws.Range("J80").Value = 14334629
$ws.Range("K80").Value = 7952.000100000001
$ws.Range("L80").Value = 43003887
$ws.Range("M80").Value = -6954.000100000001
$ws.Range("N80").Value = -43005883

$ws.Range("H83").Value = 12184832
$ws.Range("I83").Value = 2650.6667
$ws.Range("J83").Value = 14334629
$ws.Range("K83").Value = 23856.0003
$ws.Range("L83").Value = 129011661
$ws.Range("M83").Value = -18864.0003
$ws.Range("N83").Value = -129021645

$ws.Range("H86").Value = 13792.75
$ws.Range("I86").Value = 1580
$ws.Range("J86").Value = 34147.332
$ws.Range("K86").Value = 1580
$ws.Range("L86").Value = 34147.332
$ws.Range("M86").Value = -457
$ws.Range("N86").Value = -36393.332

$ws.Range("H89").Value = 13792.75
$ws.Range("I89").Value = 1580
$ws.Range("J89").Value = 34147.332
$ws.Range("K89").Value = 7900
$ws.Range("L89").Value = 170736.66
$ws.Range("M89").Value = -2284
$ws.Range("N89").Value = -181968.66

$ws.Range("H100").Value = 55557336
$ws.Range("I100").Value = 90910550
$ws.Range("J100").Value = 2286.4285
$ws.Range("K100").Value = 90910550
$ws.Range("L100").Value = 2286.4285
$ws.Range("M100").Value = -90910009
$ws.Range("N100").Value = -3368.4285

$ws.Range("H106").Value = 11497402
$ws.Range("I106").Value = 41668570
$ws.Range("K106").Value = 41668570
$ws.Range("M106").Value = -41667939

$ws.Range("H116").Value = 5070.357
$ws.Range("I116").Value = 2764.8333
$ws.Range("K116").Value = 2764.8333
$ws.Range("M116").Value = 677.1667000000002

$ws.Range("H121").Value = 17114.285
$ws.Range("I121").Value = 800
$ws.Range("J121").Value = 23640
$ws.Range("K121").Value = 2400
$ws.Range("L121").Value = 70920
$ws.Range("M121").Value = -653
$ws.Range("N121").Value = -74414

$ws.Range("H129").Value = 154769.14
$ws.Range("I129").Value = 217.6
$ws.Range("J129").Value = 167648.44
$ws.Range("K129").Value = 652.8
$ws.Range("L129").Value = 502945.32
$ws.Range("M129").Value = 4347.2
$ws.Range("N129").Value = -512945.32

$ws.Range("H132").Value = 2458.1
$ws.Range("I132").Value = 2719.4285
$ws.Range("J132").Value = 628.8
$ws.Range("K132").Value = 8158.2855
$ws.Range("L132").Value = 1886.4
$ws.Range("M132").Value = -5628.2855
$ws.Range("N132").Value = -6946.4

$ws.Range("H135").Value = 18519634
$ws.Range("I135").Value = 836.9524
$ws.Range("J135").Value = 83335420
$ws.Range("K135").Value = 7532.5716
$ws.Range("L135").Value = 750018780
$ws.Range("M135").Value = -4997.5716
$ws.Range("N135").Value = -750023850

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7701.7744
$ws.Range("I32").Value = 5474.5107
$ws.Range("K32").Value = 5474.5107
$ws.Range("M32").Value = -5187.5107

$ws.Range("H45").Value = 2319.1667
$ws.Range("I45").Value = 1612.25
$ws.Range("J45").Value = 3733
$ws.Range("K45").Value = 1612.25
$ws.Range("L45").Value = 3733
$ws.Range("M45").Value = -1235.25
$ws.Range("N45").Value = -4487

$ws.Range("H63").Value = 2919
$ws.Range("J63").Value = 2888
$ws.Range("L63").Value = 2888
$ws.Range("N63").Value = -4260

$ws.Range("H66").Value = 2919
$ws.Range("J66").Value = 2888
$ws.Range("L66").Value = 14440
$ws.Range("N66").Value = -21304

$ws.Range("H110").Value = 763.5833
$ws.Range("I110").Value = 662.6667
$ws.Range("J110").Value = 1066.3334
$ws.Range("K110").Value = 662.6667
$ws.Range("L110").Value = 1066.3334
$ws.Range("M110").Value = 1382.3333
$ws.Range("N110").Value = -5156.3334

$ws.Range("H132").Value = 10170.948
$ws.Range("I132").Value = 1470.4468
$ws.Range("J132").Value = 47345.816
$ws.Range("K132").Value = 4411.3404
$ws.Range("L132").Value = 142037.448
$ws.Range("M132").Value = -1881.3404
$ws.Range("N132").Value = -147097.448

$ws.Range("H134").Value = 62000
$ws.Range("J134").Value = 62000
$ws.Range("L134").Value = 62000
$ws.Range("N134").Value = -72140

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1730.3784
$ws.Range("I86").Value = 1616
$ws.Range("J86").Value = 1851.1111
$ws.Range("K86").Value = 1616
$ws.Range("L86").Value = 1851.1111
$ws.Range("M86").Value = -493
$ws.Range("N86").Value = -4097.1111

$ws.Range("H89").Value = 1730.3784
$ws.Range("I89").Value = 1616
$ws.Range("J89").Value = 1851.1111
$ws.Range("K89").Value = 8080
$ws.Range("L89").Value = 9255.5555
$ws.Range("M89").Value = -2464
$ws.Range("N89").Value = -20487.5555

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 327.5
$ws.Range("I22").Value = 155
$ws.Range("K22").Value = 155
$ws.Range("M22").Value = 195

$ws.Range("H86").Value = 20852968
$ws.Range("I86").Value = 3801.4
$ws.Range("K86").Value = 3801.4
$ws.Range("M86").Value = -2678.4

$ws.Range("H89").Value = 20852968
$ws.Range("I89").Value = 3801.4
$ws.Range("K89").Value = 19007
$ws.Range("M89").Value = -13391

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H130").Value = 1260
$ws.Range("I130").Value = 1260
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 3780
$ws.Range("L130").Value = 0
$ws.Range("M130").Value = 1240
$ws.Range("N130").ClearContents()

$ws.Range("H131").Value = 741.04
$ws.Range("J131").Value = 741.04
$ws.Range("L131").Value = 2223.12
$ws.Range("N131").Value = -12303.12

$ws.Range("H137").Value = 20839318
$ws.Range("I137").Value = 1851.6666
$ws.Range("J137").Value = 33341796
$ws.Range("K137").Value = 5554.9998
$ws.Range("L137").Value = 100025388
$ws.Range("M137").Value = -454.9997999999996
$ws.Range("N137").Value = -100035588

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 12954.363
$ws.Range("I70").Value = 4997.5
$ws.Range("J70").Value = 17501.143
$ws.Range("K70").Value = 4997.5
$ws.Range("L70").Value = 17501.143
$ws.Range("M70").Value = -4727.5
$ws.Range("N70").Value = -18041.143

$ws.Range("H73").Value = 12954.363
$ws.Range("I73").Value = 4997.5
$ws.Range("J73").Value = 17501.143
$ws.Range("K73").Value = 4997.5
$ws.Range("L73").Value = 17501.143
$ws.Range("M73").Value = -4061.5
$ws.Range("N73").Value = -19373.143

$ws.Range("H80").Value = 3621.611
$ws.Range("I80").Value = 3242.8572
$ws.Range("J80").Value = 3862.6365
$ws.Range("K80").Value = 3242.8572
$ws.Range("L80").Value = 3862.6365
$ws.Range("M80").Value = -2244.8572
$ws.Range("N80").Value = -5858.636500000001

$ws.Range("H83").Value = 3621.611
$ws.Range("I83").Value = 3242.8572
$ws.Range("J83").Value = 3862.6365
$ws.Range("K83").Value = 16214.286
$ws.Range("L83").Value = 19313.1825
$ws.Range("M83").Value = -11222.286
$ws.Range("N83").Value = -29297.1825

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2361.05
$ws.Range("I22").Value = 3313.923
$ws.Range("J22").Value = 591.4286
$ws.Range("K22").Value = 3313.923
$ws.Range("L22").Value = 591.4286
$ws.Range("M22").Value = -3018.923
$ws.Range("N22").Value = -1181.4286

$ws.Range("H27").Value = 2361.05
$ws.Range("I27").Value = 3313.923
$ws.Range("J27").Value = 591.4286
$ws.Range("K27").Value = 3313.923
$ws.Range("L27").Value = 591.4286
$ws.Range("M27").Value = -3206.923
$ws.Range("N27").Value = -805.4286

$ws.Range("H46").Value = 945.03845
$ws.Range("I46").Value = 962
$ws.Range("J46").Value = 815
$ws.Range("K46").Value = 962
$ws.Range("L46").Value = 815
$ws.Range("M46").Value = -774
$ws.Range("N46").Value = -1191

$ws.Range("H61").Value = 7763.3335
$ws.Range("I61").Value = 2895
$ws.Range("J61").Value = 17500
$ws.Range("K61").Value = 2895
$ws.Range("L61").Value = 17500
$ws.Range("M61").Value = -2693
$ws.Range("N61").Value = -17904

$ws.Range("H68").Value = 2998.75
$ws.Range("J68").Value = 2998.75
$ws.Range("L68").Value = 2998.75
$ws.Range("N68").Value = -4496.75

$ws.Range("H71").Value = 2998.75
$ws.Range("J71").Value = 2998.75
$ws.Range("L71").Value = 14993.75
$ws.Range("N71").Value = -22481.75

$ws.Range("H113").Value = 7763.3335
$ws.Range("I113").Value = 2895
$ws.Range("J113").Value = 17500
$ws.Range("K113").Value = 2895
$ws.Range("L113").Value = 17500
$ws.Range("M113").Value = -725
$ws.Range("N113").Value = -21840

$ws.Range("H132").Value = 310757.53
$ws.Range("I132").Value = 416698.8
$ws.Range("J132").Value = 3527.8
$ws.Range("K132").Value = 1250096.4
$ws.Range("L132").Value = 10583.4
$ws.Range("M132").Value = -1247566.4
$ws.Range("N132").Value = -15643.4

$ws.Range("H136").Value = 1735.3636
$ws.Range("I136").Value = 1638.125
$ws.Range("J136").Value = 1994.6666
$ws.Range("K136").Value = 4914.375
$ws.Range("L136").Value = 5983.9998
$ws.Range("M136").Value = -2364.375
$ws.Range("N136").Value = -11083.9998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 10083
$ws.Range("I64").Value = 8817.817999999999
$ws.Range("J64").Value = 24000
$ws.Range("K64").Value = 8817.817999999999
$ws.Range("L64").Value = 24000
$ws.Range("M64").Value = -8569.817999999999
$ws.Range("N64").Value = -24496

$ws.Range("H67").Value = 10083
$ws.Range("I67").Value = 8817.817999999999
$ws.Range("J67").Value = 24000
$ws.Range("K67").Value = 8817.817999999999
$ws.Range("L67").Value = 24000
$ws.Range("M67").Value = -7959.817999999999
$ws.Range("N67").Value = -25716

$ws.Range("H107").Value = 55024536
$ws.Range("J107").Value = 5051569.5
$ws.Range("L107").Value = 15154708.5
$ws.Range("N107").Value = -15158548.5

$ws.Range("H136").Value = 27167614
$ws.Range("I136").Value = 33299930
$ws.Range("J136").Value = 10215
$ws.Range("K136").Value = 99899790
$ws.Range("L136").Value = 30645
$ws.Range("M136").Value = -99897240
$ws.Range("N136").Value = -35745
